$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-cell value updates (Price column D / Volume(1h) column E) ---
# Each entry: row number, column letter, new value
$updates = @(
    @{ Row = 2;  Col = "D"; Value = "67.112.99" },

    @{ Row = 3;  Col = "D"; Value = "2.468.68" },
    @{ Row = 3;  Col = "E"; Value = "  +0.16%  " },

    @{ Row = 4;  Col = "E"; Value = "  -0.03%  " },

    @{ Row = 5;  Col = "D"; Value = "582.41" },
    @{ Row = 5;  Col = "E"; Value = "  -0.16%  " },

    @{ Row = 6;  Col = "D"; Value = "174.70" },
    @{ Row = 6;  Col = "E"; Value = "  +3.56%  " },

    @{ Row = 7;  Col = "E"; Value = "  -0.07%  " },

    @{ Row = 8;  Col = "D"; Value = "0.512" },
    @{ Row = 8;  Col = "E"; Value = "  -0.38%  " },

    @{ Row = 9;  Col = "D"; Value = "0.138" },
    @{ Row = 9;  Col = "E"; Value = "  +2.27%  " },

    @{ Row = 10; Col = "E"; Value = "  +0.30%  " },

    @{ Row = 11; Col = "D"; Value = "4.94" },
    @{ Row = 11; Col = "E"; Value = "  +0.25%  " },

    @{ Row = 12; Col = "D"; Value = "0.334" },
    @{ Row = 12; Col = "E"; Value = "  +1.29%  " },

    @{ Row = 13; Col = "D"; Value = "2.917.07" },
    @{ Row = 13; Col = "E"; Value = "  +0.00%  " },

    @{ Row = 14; Col = "D"; Value = "25.40" },
    @{ Row = 14; Col = "E"; Value = "  -0.62%  " },

    @{ Row = 15; Col = "D"; Value = "66.928.63" },
    @{ Row = 15; Col = "E"; Value = "  +0.35%  " },

    @{ Row = 16; Col = "E"; Value = "  +0.21%  " },

    @{ Row = 17; Col = "D"; Value = "2.470.29" },
    @{ Row = 17; Col = "E"; Value = "  +0.39%  " },

    @{ Row = 18; Col = "D"; Value = "10.93" },
    @{ Row = 18; Col = "E"; Value = "  -1.53%  " },

    @{ Row = 19; Col = "D"; Value = "7.46" },
    @{ Row = 19; Col = "E"; Value = "  -1.37%  " },

    @{ Row = 20; Col = "D"; Value = "348.46" },
    @{ Row = 20; Col = "E"; Value = "  -1.44%  " },

    @{ Row = 21; Col = "D"; Value = "3.99" },
    @{ Row = 21; Col = "E"; Value = "  -1.00%  " },

    @{ Row = 22; Col = "E"; Value = "  +0.05%  " },

    @{ Row = 23; Col = "D"; Value = "69.32" },
    @{ Row = 23; Col = "E"; Value = "  +0.54%  " },

    @{ Row = 24; Col = "D"; Value = "4.19" },
    @{ Row = 24; Col = "E"; Value = "  -1.14%  " },

    @{ Row = 25; Col = "D"; Value = "1.80" },
    @{ Row = 25; Col = "E"; Value = "  +0.28%  " },

    @{ Row = 26; Col = "D"; Value = "9.24" },
    @{ Row = 26; Col = "E"; Value = "  -0.05%  " },

    @{ Row = 27; Col = "D"; Value = "2.594.98" },
    @{ Row = 27; Col = "E"; Value = "  +0.36%  " },

    @{ Row = 28; Col = "E"; Value = "  +0.36%  " },

    @{ Row = 29; Col = "D"; Value = "0.0₃0901" },
    @{ Row = 29; Col = "E"; Value = "  -0.08%  " },

    @{ Row = 30; Col = "D"; Value = "498.98" },
    @{ Row = 30; Col = "E"; Value = "  -3.21%  " },

    @{ Row = 31; Col = "D"; Value = "7.73" },
    @{ Row = 31; Col = "E"; Value = "  -0.20%  " },

    @{ Row = 32; Col = "E"; Value = "  -0.34%  " },

    @{ Row = 33; Col = "E"; Value = "  -0.93%  " },

    @{ Row = 34; Col = "E"; Value = "  +0.02%  " },

    @{ Row = 35; Col = "E"; Value = "  +2.36%  " },

    @{ Row = 36; Col = "D"; Value = "161.41" },
    @{ Row = 36; Col = "E"; Value = "  +1.75%  " },

    @{ Row = 37; Col = "D"; Value = "18.68" },
    @{ Row = 37; Col = "E"; Value = "  +0.06%  " },

    @{ Row = 38; Col = "D"; Value = "18.15" },
    @{ Row = 38; Col = "E"; Value = "  -0.96%  " },

    @{ Row = 39; Col = "D"; Value = "1.33" },
    @{ Row = 39; Col = "E"; Value = "  -1.40%  " },

    @{ Row = 40; Col = "E"; Value = "  -0.01%  " },

    @{ Row = 41; Col = "E"; Value = "  +0.95%  " },

    @{ Row = 42; Col = "D"; Value = "0.327" },
    @{ Row = 42; Col = "E"; Value = "  +0.20%  " },

    @{ Row = 43; Col = "D"; Value = "4.82" },
    @{ Row = 43; Col = "E"; Value = "  +0.16%  " },

    @{ Row = 44; Col = "D"; Value = "2.39" },
    @{ Row = 44; Col = "E"; Value = "  +0.60%  " },

    @{ Row = 45; Col = "D"; Value = "142.57" },
    @{ Row = 45; Col = "E"; Value = "  +1.54%  " },

    @{ Row = 46; Col = "D"; Value = "3.48" },
    @{ Row = 46; Col = "E"; Value = "  +0.80%  " },

    @{ Row = 49; Col = "D"; Value = "0.0741" },
    @{ Row = 49; Col = "E"; Value = "  +1.37%  " },

    @{ Row = 50; Col = "D"; Value = "1.57" },
    @{ Row = 50; Col = "E"; Value = "  -1.14%  " },

    @{ Row = 51; Col = "D"; Value = "0.581" },
    @{ Row = 51; Col = "E"; Value = "  +0.19%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Col + $u.Row)
    # Price/volume values must stay plain text (as in the source data), so force
    # the Text number format BEFORE writing the value to stop Excel from
    # auto-coercing numeric-looking strings (e.g. "582.41") into real numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

# --- Rows 47/48 swap: BabyDogeCoin moves above ARBITRUM, both with refreshed data ---
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0255"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.509"
$ws.Range("E48").Value = "  -0.87%  "
